# Apply changes described by the diff:
#  - On the "VTQaZ" sheet, row 6 (LPG vehicle), columns M:P (years 2031-2034)
#    change value from 0 to 1.
#  - Update the selection on the "VTQaZ" sheet from Q6 to Q6:AF6.
#  - Make "About" the active sheet (tabSelected moves from VTQaZ to About).

$wb = $excel.ActiveWorkbook

$wsVTQaZ = $wb.Worksheets.Item("VTQaZ")
$wsAbout = $wb.Worksheets.Item("About")

# Update the LPG vehicle (row 6) values for years 2031-2034 (columns M-P) from 0 to 1
$wsVTQaZ.Range("M6:P6").Value = 1

# Activate VTQaZ first so we can set its selection to Q6:AF6
$wsVTQaZ.Activate()
$wsVTQaZ.Range("Q6:AF6").Select()

# Finally, make "About" the active/selected sheet
$wsAbout.Activate()
